# Auto-generated Excel COM-interop script to apply diff changes
# Updates Step1_Data, Step2_Sj (signal-value columns K..AQ for rows 2,4,6)
# and Step3_DataPts_* (columns D, F, G for rows 2,4,6) to match target workbook state.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Step1_Data")
$sheet1Updates = @(
    @("K2", 0),
    @("L2", 0.124644175938417),
    @("M2", 0.03175503718722875),
    @("N2", 0.115992988497951),
    @("O2", 0.01281427733980562),
    @("P2", 0.0000906159350017947),
    @("Q2", 0.00793182845378429),
    @("R2", 0.03859515947579245),
    @("S2", 0.06217692757305072),
    @("T2", 0.0547622795171781),
    @("U2", 0.03240954073637561),
    @("V2", 0.069972778244315),
    @("W2", 0.003837619405165615),
    @("X2", 0.08247346857849575),
    @("Y2", 0.001372999877853884),
    @("Z2", 0.08813108842243823),
    @("AA2", 0.03012969346516808),
    @("AB2", 0.003699212654417156),
    @("AC2", 0.0001621146318871093),
    @("AD2", 0.0002910323879359048),
    @("AE2", 0.02245945312733212),
    @("AF2", 0.001297558410471997),
    @("AG2", 0.1048560073874351),
    @("AH2", 0.02301759644755914),
    @("AI2", 0.01584473148031404),
    @("AJ2", 0.008160599321728393),
    @("AK2", 0.01413102225238605),
    @("AL2", 0.001404681866852471),
    @("AM2", 0.02945888801330913),
    @("AN2", 0.0032358583924201),
    @("AO2", 0.00590264526006402),
    @("AP2", 0.00002397365672861122),
    @("AQ2", 0.008964146061136682),
    @("K4", 0),
    @("L4", 0.09831669277681253),
    @("M4", 0.01582505180635173),
    @("N4", 0.1271825372158307),
    @("O4", 0.01847668261958114),
    @("P4", 0.00480704665960033),
    @("Q4", 0.009034171156222539),
    @("R4", 0.01961059351158119),
    @("S4", 0.01006267943408912),
    @("T4", 0.07937592776587296),
    @("U4", 0.006285788134006368),
    @("V4", 0.1038257779480492),
    @("W4", 0.0135677053239889),
    @("X4", 0.01807324942842082),
    @("Y4", 0.02259531884100715),
    @("Z4", 0.04917605904081197),
    @("AA4", 0.07894025272912569),
    @("AB4", 0.0001533490233310932),
    @("AC4", 0.00589247977717051),
    @("AD4", 0.01198778462468718),
    @("AE4", 0.00373473900177227),
    @("AF4", 0.002631095780442792),
    @("AG4", 0.09741563041722294),
    @("AH4", 0.07153201329593299),
    @("AI4", 0.007778395007325622),
    @("AJ4", 0.02078102312625693),
    @("AK4", 0.01516479089127121),
    @("AL4", 0.0007735251534683633),
    @("AM4", 0.05032554737315231),
    @("AN4", 0.01206517853360788),
    @("AO4", 0.009054788643471047),
    @("AP4", 0.000003824326133118951),
    @("AQ4", 0.0155503006334015),
    @("K6", 0),
    @("L6", 0.09248474330721594),
    @("M6", 0.03154883807283564),
    @("N6", 0.122077972748513),
    @("O6", 0.00479491625187038),
    @("P6", 0.002698331310644468),
    @("Q6", 0.006062074777695521),
    @("R6", 0.01375573151840904),
    @("S6", 0.008693404344043964),
    @("T6", 0.06568394666469315),
    @("U6", 0.008303326193298478),
    @("V6", 0.1216065550583814),
    @("W6", 0.01200891329975928),
    @("X6", 0.04991831341829091),
    @("Y6", 0.02646735889500729),
    @("Z6", 0.06232456796571941),
    @("AA6", 0.06446744141857676),
    @("AB6", 0.003516180208987975),
    @("AC6", 0.004575549965474717),
    @("AD6", 0.01142033550148097),
    @("AE6", 0.00721306347003378),
    @("AF6", 0.002814820063658006),
    @("AG6", 0.09614367307751866),
    @("AH6", 0.06729786993540186),
    @("AI6", 0.009851109097728987),
    @("AJ6", 0.01827483302296147),
    @("AK6", 0.01336156617769277),
    @("AL6", 0.00003650561676278285),
    @("AM6", 0.0455917093776613),
    @("AN6", 0.006707540141836981),
    @("AO6", 0.007785908200292907),
    @("AP6", 0.00004606492837252676),
    @("AQ6", 0.01246683596917956)
)
foreach ($pair in $sheet1Updates) {
    $sheet1.Range($pair[0]).Value = $pair[1]
}

$sheet2 = $wb.Worksheets.Item("Step2_Sj")
$sheet2Updates = @(
    @("K2", 0),
    @("L2", 0.124644175938417),
    @("M2", 0.1563992131256457),
    @("N2", 0.2723922016235967),
    @("O2", 0.2852064789634023),
    @("P2", 0.2852970948984041),
    @("Q2", 0.2932289233521884),
    @("R2", 0.3318240828279809),
    @("S2", 0.3940010104010316),
    @("T2", 0.4487632899182097),
    @("U2", 0.4811728306545853),
    @("V2", 0.5511456088989003),
    @("W2", 0.5549832283040659),
    @("X2", 0.6374566968825617),
    @("Y2", 0.6388296967604156),
    @("Z2", 0.7269607851828538),
    @("AA2", 0.7570904786480219),
    @("AB2", 0.7607896913024391),
    @("AC2", 0.7609518059343262),
    @("AD2", 0.7612428383222621),
    @("AE2", 0.7837022914495942),
    @("AF2", 0.7849998498600662),
    @("AG2", 0.8898558572475013),
    @("AH2", 0.9128734536950605),
    @("AI2", 0.9287181851753745),
    @("AJ2", 0.9368787844971029),
    @("AK2", 0.951009806749489),
    @("AL2", 0.9524144886163415),
    @("AM2", 0.9818733766296506),
    @("AN2", 0.9851092350220707),
    @("AO2", 0.9910118802821347),
    @("AP2", 0.9910358539388633),
    @("K4", 0),
    @("L4", 0.09831669277681253),
    @("M4", 0.1141417445831643),
    @("N4", 0.2413242817989949),
    @("O4", 0.2598009644185761),
    @("P4", 0.2646080110781764),
    @("Q4", 0.273642182234399),
    @("R4", 0.2932527757459801),
    @("S4", 0.3033154551800692),
    @("T4", 0.3826913829459422),
    @("U4", 0.3889771710799486),
    @("V4", 0.4928029490279978),
    @("W4", 0.5063706543519867),
    @("X4", 0.5244439037804075),
    @("Y4", 0.5470392226214147),
    @("Z4", 0.5962152816622267),
    @("AA4", 0.6751555343913525),
    @("AB4", 0.6753088834146835),
    @("AC4", 0.6812013631918541),
    @("AD4", 0.6931891478165413),
    @("AE4", 0.6969238868183135),
    @("AF4", 0.6995549825987563),
    @("AG4", 0.7969706130159792),
    @("AH4", 0.8685026263119121),
    @("AI4", 0.8762810213192377),
    @("AJ4", 0.8970620444454946),
    @("AK4", 0.9122268353367659),
    @("AL4", 0.9130003604902343),
    @("AM4", 0.9633259078633866),
    @("AN4", 0.9753910863969945),
    @("AO4", 0.9844458750404655),
    @("AP4", 0.9844496993665987),
    @("K6", 0),
    @("L6", 0.09248474330721594),
    @("M6", 0.1240335813800516),
    @("N6", 0.2461115541285646),
    @("O6", 0.250906470380435),
    @("P6", 0.2536048016910794),
    @("Q6", 0.259666876468775),
    @("R6", 0.273422607987184),
    @("S6", 0.282116012331228),
    @("T6", 0.3477999589959211),
    @("U6", 0.3561032851892196),
    @("V6", 0.477709840247601),
    @("W6", 0.4897187535473603),
    @("X6", 0.5396370669656512),
    @("Y6", 0.5661044258606585),
    @("Z6", 0.6284289938263778),
    @("AA6", 0.6928964352449546),
    @("AB6", 0.6964126154539426),
    @("AC6", 0.7009881654194173),
    @("AD6", 0.7124085009208982),
    @("AE6", 0.7196215643909321),
    @("AF6", 0.72243638445459),
    @("AG6", 0.8185800575321087),
    @("AH6", 0.8858779274675106),
    @("AI6", 0.8957290365652396),
    @("AJ6", 0.9140038695882011),
    @("AK6", 0.9273654357658939),
    @("AL6", 0.9274019413826567),
    @("AM6", 0.9729936507603181),
    @("AN6", 0.979701190902155),
    @("AO6", 0.9874870991024479),
    @("AP6", 0.9875331640308205)
)
foreach ($pair in $sheet2Updates) {
    $sheet2.Range($pair[0]).Value = $pair[1]
}

$sheet3 = $wb.Worksheets.Item("Step3_DataPts_0.5")
$sheet3Updates = @(
    @("F2", 0.5511456088989003),
    @("D4", 22),
    @("F4", 0.5063706543519867),
    @("G4", 13),
    @("F6", 0.5396370669656512)
)
foreach ($pair in $sheet3Updates) {
    $sheet3.Range($pair[0]).Value = $pair[1]
}

$sheet4 = $wb.Worksheets.Item("Step3_DataPts_0.7")
$sheet4Updates = @(
    @("F2", 0.7269607851828538),
    @("D4", 32),
    @("F4", 0.7969706130159792),
    @("G4", 23),
    @("D6", 28),
    @("F6", 0.7009881654194173),
    @("G6", 19)
)
foreach ($pair in $sheet4Updates) {
    $sheet4.Range($pair[0]).Value = $pair[1]
}

$sheet5 = $wb.Worksheets.Item("Step3_DataPts_0.8")
$sheet5Updates = @(
    @("F2", 0.8898558572475013),
    @("D4", 33),
    @("F4", 0.8685026263119121),
    @("G4", 24),
    @("F6", 0.8185800575321087)
)
foreach ($pair in $sheet5Updates) {
    $sheet5.Range($pair[0]).Value = $pair[1]
}

$sheet6 = $wb.Worksheets.Item("Step3_DataPts_0.9")
$sheet6Updates = @(
    @("F2", 0.9128734536950605),
    @("D4", 36),
    @("F4", 0.9122268353367659),
    @("G4", 27),
    @("D6", 35),
    @("F6", 0.9140038695882011),
    @("G6", 26)
)
foreach ($pair in $sheet6Updates) {
    $sheet6.Range($pair[0]).Value = $pair[1]
}

Write-Host "Applied all cell updates."